# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 and 17: swap the worker identity (doc number + name) and swap
# the "Valor Mora" (column G) between the two rows.
$ws.Range("C16").Value = "1050963433"
$ws.Range("D16").Value = "TANIA MARGARITA RODRIGUEZ BARBOZA"
$ws.Range("G16").Value = 828116

$ws.Range("C17").Value = "1143399179"
$ws.Range("D17").Value = "JULIO NICOLAS MESA ZABALETA"
$ws.Range("G17").Value = 781242

# Rows 19-23: reverse the "Periodo Mora" values (1908..1912 -> 1912..1908)
$ws.Range("E19").Value = "1912"
$ws.Range("E20").Value = "1911"
$ws.Range("E21").Value = "1910"
$ws.Range("E22").Value = "1909"
$ws.Range("E23").Value = "1908"

# Columns auto-fit to the new data widths (bestFit columns B, C, E, F, G, H, I, J)
$ws.Range("B15:J24").Columns.AutoFit()
